$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying the formatting (style) of the existing
# header cell G1 so the new column matches the rest of the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill the new "Save" column (H2:H8) with 1 for every data row.
$ws.Range("H2:H8").Value = 1
